$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows appended below the existing "Stephon Gilmore" block.
# Rows 5-7:  Denzel Ward   -> new highlight colour (light yellow)
# Rows 8-10: Isaac Yiadom  -> same highlight colour already used by the
#                             "Stephon Gilmore" rows (light green)
# ---------------------------------------------------------------------------

$data = @(
    @(5,  "Denzel Ward",  "Group1",     13,                 44.33333333333334, 36.66666666666666, 7.666666666666667),
    @(6,  "Denzel Ward",  "Group2",     15,                 45.33333333333334, 37.66666666666666, 7.666666666666667),
    @(7,  "Denzel Ward",  "Difference", 2,                  1,                 1,                  0),
    @(8,  "Isaac Yiadom", "Group1",     3,                  33,                23.33333333333333, 9.666666666666666),
    @(9,  "Isaac Yiadom", "Group2",     6.666666666666667,  28.11111111111111, 21.22222222222222, 6.888888888888888),
    @(10, "Isaac Yiadom", "Difference", 3.666666666666667, -4.888888888888886,-2.111111111111107,-2.777777777777778)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

# Highlight the new "Denzel Ward" rows with a new light-yellow fill.
$ws.Range("A5:F7").Interior.Color = 12451839

# Re-use the existing light-green fill (already applied to the
# "Stephon Gilmore" rows) for the "Isaac Yiadom" rows by copying its
# formatting across.
$ws.Range("A2:F2").Copy()
$ws.Range("A8:F10").PasteSpecial(-4122)
